$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update row 2 odds that changed ---
$ws.Range("I2").Value = 6.5
$ws.Range("AC2").Value = 21
$ws.Range("AH2").Value = 21
$ws.Range("BD2").Value = 151

# --- Row 9 (Poland Ekstraklasa) is replaced with what used to be row 11 (Scotland - Championship) ---
$row9 = @(
    'llm6eDM8',
    '22/11/2024',
    '16:45',
    'SCOTLAND - CHAMPIONSHIP',
    'Queen''s Park',
    'Falkirk',
    3.9,
    4.1,
    1.73,
    4.33,
    2.3,
    2.25,
    1.03,
    10,
    1.2,
    4.33,
    1.62,
    2.25,
    1.3,
    3.4,
    1.67,
    2.1,
    15,
    23,
    13,
    41,
    29,
    34,
    15,
    8.5,
    15,
    41,
    151,
    9,
    9.5,
    8.5,
    15,
    13,
    21,
    6.5,
    21,
    26,
    67,
    81,
    151,
    3.4,
    7.5,
    41,
    4,
    9,
    17,
    26,
    41,
    101,
    81,
    81
)
for ($i = 0; $i -lt $row9.Length; $i++) {
    $ws.Cells.Item(9, $i + 1).Value = $row9[$i]
}

# --- Row 10 (Poland Division 1) is replaced with what used to be row 13 (Wales - Cymru Premier) ---
$row10 = @(
    'ABXWOPog',
    '22/11/2024',
    '16:45',
    'WALES - CYMRU PREMIER',
    'Newtown',
    'Connahs Q.',
    2.52,
    3.5,
    2.4,
    3.1,
    2.22,
    2.95,
    1.04,
    8.25,
    1.24,
    3.7,
    1.72,
    2.05,
    1.34,
    3,
    1.62,
    2.18,
    10,
    14,
    9.75,
    28,
    19.5,
    26,
    8.25,
    7,
    13,
    50,
    350,
    9.75,
    13,
    9.5,
    25,
    18.5,
    25,
    4.65,
    13,
    19.5,
    55,
    80,
    200,
    3,
    6.8,
    55,
    4.5,
    12.5,
    19,
    50,
    75,
    200,
    '',
    ''
)
for ($i = 0; $i -lt $row10.Length; $i++) {
    $ws.Cells.Item(10, $i + 1).Value = $row10[$i]
}

# --- Remove the now-duplicated trailing rows (old rows 11, 12, 13) ---
$ws.Rows.Item(13).Delete()
$ws.Rows.Item(12).Delete()
$ws.Rows.Item(11).Delete()

Write-Host "Edit complete"
